# Apply the cryptos.xlsx price/volume update (commit: "Updated cryptos list on
# Sun May 12 04:56:32 UTC 2024 with GitHub Actions").
#
# All target cells are plain (non-formula) strings (t="inlineStr" in the source
# OOXML) with no cell-style changes. Column D ("Price") sometimes holds values
# that LOOK like plain numbers (e.g. "33.69"); a bare `Range.Value = ...` on such
# a string lets Excel auto-convert it to a real number (and, via the quote-prefix
# bookkeeping, tags the cell with a new "@ text" style). To keep these cells
# text-typed with their ORIGINAL (unstyled) formatting we:
#   1. force the cell to Text format before writing ("@"),
#   2. write the literal string,
#   3. reset the cell style back to "Normal" (style index 0, matching the source).
# Cells whose new text can never be mis-parsed as a number (it has two dots,
# letters, a "%", a URL, etc.) are just assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.055.25"
$ws.Range("D3").Value = "2.927.21"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.57%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("E9").Value = "  +1.61%  "
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("E11").Value = "  -1.13%  "
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").Value = "3.411.25"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "60.937.69"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.926.26"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "432.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("E20").Value = "  -1.56%  "
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.67"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("E32").Value = "  +2.77%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").Value = "0.0$([char]0x2083)0861"
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("E39").Value = "  -4.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.283"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "379.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.97%  "
$ws.Range("D44").Value = "2.703.26"
$ws.Range("E44").Value = "  +1.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0343"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("E50").Value = "  -2.35%  "
$ws.Range("E51").Value = "  -0.56%  "
